# Updates the EC (Estado de Cuenta) data: reorders rows 16-36 from
# "grouped by worker" (7 periods each for 3 workers) to
# "grouped by period" (3 workers each for 7 periods, ascending 2302->2308).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "1051419192"
$ws.Range("D16").Value = "CARLOS ALFONSO CASTILLO PAJARO"
$ws.Range("E16").Value = "2302"
$ws.Range("C17").Value = "1235038902"
$ws.Range("D17").Value = "SERGIO JOSE CABALLERO OSPINO"
$ws.Range("E17").Value = "2302"
$ws.Range("C18").Value = "1003344794"
$ws.Range("D18").Value = "ANDRES FELIPE BARRIOS ECHEVERRIA"
$ws.Range("E18").Value = "2302"
$ws.Range("C19").Value = "1051419192"
$ws.Range("D19").Value = "CARLOS ALFONSO CASTILLO PAJARO"
$ws.Range("E19").Value = "2303"
$ws.Range("C20").Value = "1235038902"
$ws.Range("D20").Value = "SERGIO JOSE CABALLERO OSPINO"
$ws.Range("E20").Value = "2303"
$ws.Range("C21").Value = "1003344794"
$ws.Range("D21").Value = "ANDRES FELIPE BARRIOS ECHEVERRIA"
$ws.Range("E21").Value = "2303"
$ws.Range("C22").Value = "1051419192"
$ws.Range("D22").Value = "CARLOS ALFONSO CASTILLO PAJARO"
$ws.Range("E22").Value = "2304"
$ws.Range("C23").Value = "1235038902"
$ws.Range("D23").Value = "SERGIO JOSE CABALLERO OSPINO"
$ws.Range("E23").Value = "2304"
$ws.Range("C24").Value = "1003344794"
$ws.Range("D24").Value = "ANDRES FELIPE BARRIOS ECHEVERRIA"
$ws.Range("E24").Value = "2304"
$ws.Range("C25").Value = "1051419192"
$ws.Range("D25").Value = "CARLOS ALFONSO CASTILLO PAJARO"
$ws.Range("E25").Value = "2305"
$ws.Range("C26").Value = "1235038902"
$ws.Range("D26").Value = "SERGIO JOSE CABALLERO OSPINO"
$ws.Range("E26").Value = "2305"
$ws.Range("C27").Value = "1003344794"
$ws.Range("D27").Value = "ANDRES FELIPE BARRIOS ECHEVERRIA"
$ws.Range("E27").Value = "2305"
$ws.Range("C28").Value = "1051419192"
$ws.Range("D28").Value = "CARLOS ALFONSO CASTILLO PAJARO"
$ws.Range("E28").Value = "2306"
$ws.Range("C29").Value = "1235038902"
$ws.Range("D29").Value = "SERGIO JOSE CABALLERO OSPINO"
$ws.Range("E29").Value = "2306"
$ws.Range("C30").Value = "1003344794"
$ws.Range("D30").Value = "ANDRES FELIPE BARRIOS ECHEVERRIA"
$ws.Range("E30").Value = "2306"
$ws.Range("C31").Value = "1051419192"
$ws.Range("D31").Value = "CARLOS ALFONSO CASTILLO PAJARO"
$ws.Range("E31").Value = "2307"
$ws.Range("C32").Value = "1235038902"
$ws.Range("D32").Value = "SERGIO JOSE CABALLERO OSPINO"
$ws.Range("E32").Value = "2307"
$ws.Range("C33").Value = "1003344794"
$ws.Range("D33").Value = "ANDRES FELIPE BARRIOS ECHEVERRIA"
$ws.Range("E33").Value = "2307"
$ws.Range("C34").Value = "1051419192"
$ws.Range("D34").Value = "CARLOS ALFONSO CASTILLO PAJARO"
$ws.Range("E34").Value = "2308"
$ws.Range("C35").Value = "1235038902"
$ws.Range("D35").Value = "SERGIO JOSE CABALLERO OSPINO"
$ws.Range("E35").Value = "2308"
$ws.Range("C36").Value = "1003344794"
$ws.Range("D36").Value = "ANDRES FELIPE BARRIOS ECHEVERRIA"
$ws.Range("E36").Value = "2308"
